$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: row 3 (06ee1f42... entry), "Latest HO Xliff Generate Date"
$wsOverview.Range("G3").Value = "2016-09-05 12:56:20"

# zh-cn sheet: row 3 (06ee1f42... entry)
$wsZhCn.Range("H3").Value = "2016-09-05 12:56:15"
$wsZhCn.Range("K3").Value = "2016-09-05 12:56:34"

# de-de sheet: row 3 (06ee1f42... entry)
$wsDeDe.Range("H3").Value = "2016-09-05 12:56:20"
$wsDeDe.Range("K3").Value = "2016-09-05 12:56:42"
